# Auto-generated edit script: update crypto price/volume columns (D, E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "58.462.76"
$ws.Cells.Item(2, 5).Value = "  +1.04%  "
$ws.Cells.Item(3, 4).Value = "3.148.74"
$ws.Cells.Item(3, 5).Value = "  +0.27%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Formula = "'537.01"
$ws.Cells.Item(5, 5).Value = "  +1.32%  "
$ws.Cells.Item(6, 4).Formula = "'140.32"
$ws.Cells.Item(6, 5).Value = "  +1.22%  "
$ws.Cells.Item(7, 5).Value = "  -0.12%  "
$ws.Cells.Item(8, 4).Value = "3.148.07"
$ws.Cells.Item(8, 5).Value = "  +0.25%  "
$ws.Cells.Item(9, 4).Formula = "'0.469"
$ws.Cells.Item(9, 5).Value = "  +4.52%  "
$ws.Cells.Item(10, 5).Value = "  +1.31%  "
$ws.Cells.Item(11, 5).Value = "  +0.06%  "
$ws.Cells.Item(12, 4).Formula = "'0.417"
$ws.Cells.Item(12, 5).Value = "  +4.48%  "
$ws.Cells.Item(13, 4).Value = "3.684.92"
$ws.Cells.Item(13, 5).Value = "  +0.04%  "
$ws.Cells.Item(14, 5).Value = "  +1.33%  "
$ws.Cells.Item(15, 4).Formula = "'25.88"
$ws.Cells.Item(15, 5).Value = "  +1.20%  "
$ws.Cells.Item(16, 5).Value = "  +0.43%  "
$ws.Cells.Item(17, 4).Value = "58.526.75"
$ws.Cells.Item(17, 5).Value = "  +0.90%  "
$ws.Cells.Item(18, 4).Value = "3.146.47"
$ws.Cells.Item(18, 5).Value = "  +0.05%  "
$ws.Cells.Item(19, 4).Formula = "'6.09"
$ws.Cells.Item(19, 5).Value = "  +1.56%  "
$ws.Cells.Item(20, 4).Formula = "'12.78"
$ws.Cells.Item(20, 5).Value = "  +0.22%  "
$ws.Cells.Item(21, 4).Formula = "'8.23"
$ws.Cells.Item(21, 5).Value = "  +3.06%  "
$ws.Cells.Item(22, 4).Formula = "'362.00"
$ws.Cells.Item(22, 5).Value = "  +2.60%  "
$ws.Cells.Item(23, 5).Value = "  -0.07%  "
$ws.Cells.Item(24, 4).Formula = "'69.24"
$ws.Cells.Item(24, 5).Value = "  +0.70%  "
$ws.Cells.Item(25, 4).Formula = "'0.509"
$ws.Cells.Item(25, 5).Value = "  +0.12%  "
$ws.Cells.Item(26, 5).Value = "  -1.11%  "
$ws.Cells.Item(27, 4).Formula = "'1.00"
$ws.Cells.Item(27, 5).Value = "  -0.20%  "
$ws.Cells.Item(28, 4).Value = "0.0₃0888"
$ws.Cells.Item(28, 5).Value = "  -3.59%  "
$ws.Cells.Item(29, 4).Formula = "'7.39"
$ws.Cells.Item(29, 5).Value = "  -2.19%  "
$ws.Cells.Item(30, 4).Formula = "'6.19"
$ws.Cells.Item(30, 5).Value = "  +0.22%  "
$ws.Cells.Item(31, 5).Value = "  -0.22%  "
$ws.Cells.Item(32, 4).Formula = "'21.59"
$ws.Cells.Item(32, 5).Value = "  +1.85%  "
$ws.Cells.Item(33, 4).Formula = "'5.17"
$ws.Cells.Item(33, 5).Value = "  +3.29%  "
$ws.Cells.Item(34, 5).Value = "  -2.46%  "
$ws.Cells.Item(35, 4).Formula = "'158.95"
$ws.Cells.Item(35, 5).Value = "  +0.65%  "
$ws.Cells.Item(36, 4).Formula = "'6.13"
$ws.Cells.Item(36, 5).Value = "  -0.99%  "
$ws.Cells.Item(37, 4).Formula = "'26.18"
$ws.Cells.Item(37, 5).Value = "  -1.03%  "
$ws.Cells.Item(38, 5).Value = "  +1.93%  "
$ws.Cells.Item(39, 5).Value = "  +4.73%  "
$ws.Cells.Item(40, 4).Formula = "'0.0676"
$ws.Cells.Item(40, 5).Value = "  +0.49%  "
$ws.Cells.Item(41, 4).Value = "2.501.98"
$ws.Cells.Item(41, 5).Value = "  +6.72%  "
$ws.Cells.Item(42, 4).Formula = "'0.706"
$ws.Cells.Item(42, 5).Value = "  -0.03%  "
$ws.Cells.Item(43, 4).Formula = "'4.04"
$ws.Cells.Item(43, 5).Value = "  -4.09%  "
$ws.Cells.Item(44, 4).Formula = "'37.71"
$ws.Cells.Item(44, 5).Value = "  +2.91%  "
$ws.Cells.Item(45, 4).Value = "3.189.39"
$ws.Cells.Item(45, 5).Value = "  +0.20%  "
$ws.Cells.Item(46, 4).Formula = "'0.0270"
$ws.Cells.Item(46, 5).Value = "  -1.29%  "
$ws.Cells.Item(47, 5).Value = "  -0.05%  "
$ws.Cells.Item(48, 5).Value = "  +3.38%  "
$ws.Cells.Item(49, 4).Formula = "'6.09"
$ws.Cells.Item(49, 5).Value = "  +0.72%  "
$ws.Cells.Item(50, 4).Formula = "'20.05"
$ws.Cells.Item(50, 5).Value = "  -1.81%  "
$ws.Cells.Item(51, 5).Value = "  -3.07%  "
